# Updates the NATMI Fgf9-Fgfr2 LR-pair sheet with refreshed TPM-derived
# stats. Sending/target clusters are now ECs x {ECs,FAPs,MuSCs,Resolving-Mac}
# and MuSCs x {ECs,FAPs,MuSCs,Resolving-Mac} (8 data rows instead of the
# original 4 ECs x {ECs,FAPs,Inflammatory-Mac,MuSCs} rows).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data table (rows 2-9, columns A-T)
$data = @(
    @("ECs", "Fgf9", "Fgfr2", "ECs", 2, 0.6666666666666666, 0.7364236666666667, 2.209271, 0.9917500467982164, 0.9944848646626661, 3, 1, 0.2858606666666667, 0.857582, 0.0687156860066334, 0.06932858672617494, 0.2105145603024444, 1.894631042722, 0.06814878481285022, 0.06894623018763399),
    @("ECs", "Fgf9", "Fgfr2", "FAPs", 2, 0.6666666666666666, 0.7364236666666667, 2.209271, 0.9917500467982164, 0.9944848646626661, 3, 1, 3.763360333333333, 11.290081, 0.90464312565499, 0.9127119736118995, 2.771427615661223, 24.942848540951, 0.8971798622040211, 0.9076782435534247),
    @("ECs", "Fgf9", "Fgfr2", "MuSCs", 2, 0.6666666666666666, 0.7364236666666667, 2.209271, 0.9917500467982164, 0.9944848646626661, 2, 1, 0.110331, 0.220662, 0.02652155835639462, 0.01783874265571248, 0.081250359567, 0.487502157402, 0.02630275674111599, 0.01774035957571835),
    @("ECs", "Fgf9", "Fgfr2", "Resolving-Mac", 2, 0.6666666666666666, 0.7364236666666667, 2.209271, 0.9917500467982164, 0.9944848646626661, 1, 0.3333333333333333, 0.0004976666666666667, 0.001493, 0.0001196299819817856, 0.0001206970062130259, 0.0003664935114444444, 0.003298441603, 0.0001186430402289057, 0.00012003134588895),
    @("MuSCs", "Fgf9", "Fgfr2", "ECs", 1, 0.5, 0.006126, 0.012252, 0.008249953201783585, 0.005515135337333892, 3, 1, 0.2858606666666667, 0.857582, 0.0687156860066334, 0.06932858672617494, 0.001751182444, 0.010507094664, 0.0005669011937831808, 0.0003823565385409448),
    @("MuSCs", "Fgf9", "Fgfr2", "FAPs", 1, 0.5, 0.006126, 0.012252, 0.008249953201783585, 0.005515135337333892, 3, 1, 3.763360333333333, 11.290081, 0.90464312565499, 0.9127119736118995, 0.023054345402, 0.138326072412, 0.007463263450968895, 0.005033730058474745),
    @("MuSCs", "Fgf9", "Fgfr2", "MuSCs", 1, 0.5, 0.006126, 0.012252, 0.008249953201783585, 0.005515135337333892, 2, 1, 0.110331, 0.220662, 0.02652155835639462, 0.01783874265571248, 0.000675887706, 0.002703550824, 0.000218801615278628, [double]"9.83830799941253e-05"),
    @("MuSCs", "Fgf9", "Fgfr2", "Resolving-Mac", 1, 0.5, 0.006126, 0.012252, 0.008249953201783585, 0.005515135337333892, 1, 0.3333333333333333, 0.0004976666666666667, 0.001493, 0.0001196299819817856, 0.0001206970062130259, [double]"3.048706e-06", [double]"1.8292236e-05", [double]"9.869417528799447e-07", [double]"6.656603240758675e-07")

)

$startRow = 2
for ($r = 0; $r -lt $data.Length; $r++) {
    $rowVals = $data[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $rowVals.Length; $c++) {
        $ws.Cells.Item($excelRow, $c + 1).Value = $rowVals[$c]
    }
}
